$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.395.40'
$ws.Range("E2").Value = '  +2.91%  '

$ws.Range("D3").Value = '2.306.56'
$ws.Range("E3").Value = '  +1.81%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '311.16'
$ws.Range("E5").Value = '  +1.55%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '102.43'
$ws.Range("E6").Value = '  +5.88%  '

$ws.Range("E7").Value = '  +1.49%  '

$ws.Range("E8").Value = '  -0.03%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.529'
$ws.Range("E9").Value = '  +7.51%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.76'
$ws.Range("E10").Value = '  +2.09%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0813'
$ws.Range("E11").Value = '  +2.90%  '

$ws.Range("E12").Value = '  -0.88%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.98'
$ws.Range("E13").Value = '  +0.90%  '

$ws.Range("D14").Value = '2.666.17'
$ws.Range("E14").Value = '  +1.84%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.00'
$ws.Range("E15").Value = '  +2.21%  '

$ws.Range("D16").Value = '2.313.13'
$ws.Range("E16").Value = '  +2.08%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.808'
$ws.Range("E17").Value = '  +1.99%  '

$ws.Range("D18").Value = '43.297.24'
$ws.Range("E18").Value = '  +2.95%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.33'
$ws.Range("E19").Value = '  +0.39%  '

$ws.Range("E20").Value = '  +3.11%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.19'
$ws.Range("E21").Value = '  +3.06%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '68.10'
$ws.Range("E22").Value = '  +0.46%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '241.34'
$ws.Range("E23").Value = '  +1.78%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.63'
$ws.Range("E24").Value = '  +2.11%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.01'
$ws.Range("E25").Value = '  +2.38%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.60'
$ws.Range("E27").Value = '  +4.66%  '

$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.30'
$ws.Range("E28").Value = '  +8.78%  '

$ws.Range("B29").Value = 'InjectiveProtocol'
$ws.Range("C29").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '37.00'
$ws.Range("E29").Value = '  -1.80%  '

$ws.Range("E30").Value = '  +0.75%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '168.48'
$ws.Range("E31").Value = '  +3.38%  '

$ws.Range("E33").Value = '  +0.07%  '

$ws.Range("E34").Value = '  +6.15%  '

$ws.Range("E35").Value = '  +0.74%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '17.66'
$ws.Range("E36").Value = '  -0.03%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.07'
$ws.Range("E37").Value = '  -3.77%  '

$ws.Range("E38").Value = '  +3.37%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.105'
$ws.Range("E39").Value = '  +0.97%  '

$ws.Range("E40").Value = '  +1.62%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.36'
$ws.Range("E41").Value = '  +7.42%  '

$ws.Range("E42").Value = '  -0.56%  '

$ws.Range("B43").Value = 'EnergySwap'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '19.50'
$ws.Range("E43").Value = '  +2.65%  '

$ws.Range("D44").Value = '1.973.56'
$ws.Range("E44").Value = '  +1.25%  '

$ws.Range("B45").Value = 'VeChain'
$ws.Range("C45").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0289'
$ws.Range("E45").Value = '  +2.99%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.98'
$ws.Range("E46").Value = '  +1.80%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.85'
$ws.Range("E47").Value = '  +0.07%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '55.35'
$ws.Range("E48").Value = '  +2.23%  '

$ws.Range("E49").Value = '  +1.53%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.58'
$ws.Range("E50").Value = '  +7.47%  '

$ws.Range("D51").Value = '2.534.73'
$ws.Range("E51").Value = '  +1.80%  '
